$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Move the "_GoBack" bookmark.
#
# In the original document it sits (collapsed) at the very end of the
# "1. 6 đặc tính của dữ liệu:" paragraph. In the edited document it
# instead sits in the middle of the word "ứng dụng." a few paragraphs
# below, splitting that run into "...ứng dụ" + "ng.". Word only ever
# keeps a single "_GoBack" bookmark in a document, so (re)adding one
# at the new spot automatically removes the old one - do this first,
# while character offsets for the lower paragraph are still the
# original ones (it does not change the document's character count).
# ------------------------------------------------------------------
$anchor = $d.Content
$gotAnchor = $anchor.Find.Execute("hỏi tương lai của ứng dụng", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $gotAnchor) {
    throw "Could not find the '...ứng dụng' anchor text for the _GoBack bookmark"
}
$splitPos = $anchor.End - 2   # right before the final "ng." of "ứng dụng."
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Step 2: Rewrite the bold heading.
#
# "6 đặc tính của dữ liệu:" -> "Hãy nêu và phân tích các đặc tính dữ
# liệu của yêu cầu:", restyled with the "fontstyle01" character style
# and theme-based fonts.
# ------------------------------------------------------------------
$heading = $d.Content
$gotHeading = $heading.Find.Execute("6 đặc tính của dữ liệu:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $gotHeading) {
    throw "Could not find the heading text to replace"
}
$headingRange = $d.Range($heading.Start, $heading.End)

$newHeadingXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:r><w:rPr><w:rStyle w:val="fontstyle01"/><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Hãy nêu và phân tích các đặc tính dữ liệu của yêu cầ</w:t></w:r>
<w:r><w:rPr><w:rStyle w:val="fontstyle01"/><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>u:</w:t></w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$headingRange.InsertXML($newHeadingXml)

Write-Host "Done."
